$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing "locked" (B) and "received_final" (J) values
# for rows 61-63 (rounds 49-51), matching the pattern already used in rows <=60.
$ws.Range("B61").Value = 0
$ws.Range("J61").Value = 1

$ws.Range("B62").Value = 0
$ws.Range("J62").Value = 1

$ws.Range("B63").Value = 0
$ws.Range("J63").Value = 1

# Add new row 64 for round 52 (panel E, wave 17)
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = "uk"
$ws.Range("D64").Value = 52
$ws.Range("E64").Value = "E"
$ws.Range("F64").Value = 17

# Copy the date cell formatting from the row above so we reuse the existing
# date style instead of creating a new one, then set the actual date value.
$ws.Range("G63").Copy($ws.Range("G64"))
$ws.Range("G64").Value = (Get-Date -Year 2021 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)

$ws.Range("H64").Value = "20-100562_PEW17_Final_ICUO"
$ws.Range("I64").Formula = '=C64&"_"&"sr"&TEXT(D64,"00")&"_"&YEAR(G64)&TEXT(G64,"MM")&TEXT(G64,"DD")&"_p"&E64&"_wv"&TEXT(F64,"00")&""'
$ws.Range("J64").Value = 1

# Reflect the final view state: scrolled down to show the new row, with the
# newly entered "received_final" cell selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J64").Select() | Out-Null
